# Apply weekly price/volume/date shuffle for Fruta/Vega Modelo de Temuco - Guayaba
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2..13) for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$data = @{
    2  = @{ D = 44432; M = 30;  N = 1300; O = 1300; P = 1300; S = 1300 }
    3  = @{ D = 44438; M = 60;  N = 1200; O = 1200; P = 1200; S = 1200 }
    4  = @{ D = 44424; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    5  = @{ D = 44431; M = 100; N = 1300; O = 1300; P = 1300; S = 1300 }
    6  = @{ D = 44435; M = 130; N = 1300; O = 1300; P = 1300; S = 1300 }
    7  = @{ D = 44417; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    8  = @{ D = 44405; M = 50;  N = 1200; O = 1200; P = 1200; S = 1200 }
    9  = @{ D = 44418; M = 40;  N = 1200; O = 1200; P = 1200; S = 1200 }
    10 = @{ D = 44357; M = 35;  N = 1000; O = 1000; P = 1000; S = 1000 }
    11 = @{ D = 44476; M = 80;  N = 1200; O = 1200; P = 1200; S = 1200 }
    12 = @{ D = 44343; M = 60;  N = 1300; O = 1300; P = 1300; S = 1300 }
    13 = @{ D = 44473; M = 120; N = 1200; O = 1200; P = 1200; S = 1200 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
